# Updated cryptos list values (Price and Volume(1h)) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.007.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.93%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.513.07"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.24"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.94"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.29%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.142"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.07%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.13%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.126.19"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +12.81%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.997.98"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.525.31"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.76"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "397.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.66"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.546"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.70"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.53%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.99"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.43"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.69"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.57%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.68"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.55%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.877"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.15"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.69%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.81"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.70"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.74"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0736"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.813.90"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0305"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "344.51"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.76%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.63"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.46%  "

